# Update Name of Algo
# Apply corrected values to the result_data_KNN sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value  = -12.634
$ws.Range("C3").Value  = -12.634
$ws.Range("C5").Value  = -12.634
$ws.Range("D7").Value  = -7.151999999999999
$ws.Range("A9").Value  = -20.775
$ws.Range("D9").Value  = -7.873
$ws.Range("C11").Value = -12.917
$ws.Range("C12").Value = -12.628
$ws.Range("A13").Value = -21.99
$ws.Range("A16").Value = -20.916
$ws.Range("A18").Value = -21.577
$ws.Range("A20").Value = -21.682
$ws.Range("C21").Value = -12.37
$ws.Range("D21").Value = -7.761
